$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 ("Team13") is the template row for the two new rows being appended.
# Insert two new rows (13 and 14) as copies of row 12 (values + styles),
# shifting everything below down - this reproduces the exact per-cell
# style indices used on row 12.
$ws.Rows.Item(12).Copy()
$ws.Rows.Item(13).Insert(-4121)

$ws.Rows.Item(12).Copy()
$ws.Rows.Item(14).Insert(-4121)

# Update the Team Name (column A) for the two freshly inserted rows.
$ws.Range("A13").Value = "New Team"
$ws.Range("A14").Value = "X Team13"

# Re-create the Coach Email / Member1 Email hyperlinks for the new rows,
# mirroring the ones already present on row 12.
$ws.Hyperlinks.Add($ws.Range("F13"), "mailto:jecile7288@netjook.com", "", "", "jecile7288@netjook.com")
$ws.Hyperlinks.Add($ws.Range("H13"), "mailto:safwan.du16@gmail.com", "", "", "safwan.du16@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F14"), "mailto:jecile7288@netjook.com", "", "", "jecile7288@netjook.com")
$ws.Hyperlinks.Add($ws.Range("H14"), "mailto:safwan.du16@gmail.com", "", "", "safwan.du16@gmail.com")

# Move the active selection to the last inserted row, as in the author's edit.
$ws.Range("A14").Select() | Out-Null
